$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A158").Value = "IMX-USD"
$ws.Range("A159").Value = "GRT-USD"
